$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Q8" header in column J, copying the header style (bold, centered, bordered)
# from the existing I1 header cell.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J1").Value = "Q8"

# Update the simulated/"bugfixed" rt_data values for each quarter row (A:J),
# matching the newly re-evaluated naive AVERAGE(10,9) qoq errors.
$ws.Range("A2").Value = "2022-01-01 00:00:00_diff"
$ws.Range("B2").Value = -4.324682701351703
$ws.Range("C2").Value = 3.660091258637736
$ws.Range("D2").Value = 1.731369225691766
$ws.Range("E2").Value = 3.149790629511104
$ws.Range("F2").Value = -3.635177286302406
$ws.Range("G2").Value = 0.06544417180023943
$ws.Range("H2").Value = -0.1340705272443026
$ws.Range("A3").Value = "2022-04-01 00:00:00_diff"
$ws.Range("B3").Value = 3.605598293642375
$ws.Range("C3").Value = 1.676876260696405
$ws.Range("D3").Value = 3.095297664515743
$ws.Range("E3").Value = -3.689670251297767
$ws.Range("F3").Value = 0.01095120680487852
$ws.Range("G3").Value = -0.1885634922396635
$ws.Range("A4").Value = "2022-07-01 00:00:00_diff"
$ws.Range("B4").Value = 0.8917107665847916
$ws.Range("C4").Value = 2.310132170404129
$ws.Range("D4").Value = -4.47483574540938
$ws.Range("E4").Value = -0.7742142873067352
$ws.Range("F4").Value = -0.9737289863512772
$ws.Range("G4").Value = 1.180395370885164
$ws.Range("H4").Value = -3.253940853883919
$ws.Range("I4").Value = 1.466850648880309
$ws.Range("J4").Value = -0.1948090203170254
$ws.Range("A5").Value = "2022-10-01 00:00:00_diff"
$ws.Range("B5").Value = 2.584158590043411
$ws.Range("C5").Value = -4.200809325770098
$ws.Range("D5").Value = -0.5001878676674536
$ws.Range("E5").Value = -0.6997025667119956
$ws.Range("F5").Value = 1.454421790524445
$ws.Range("G5").Value = -2.979914434244638
$ws.Range("H5").Value = 1.74087706851959
$ws.Range("I5").Value = 0.07921739932225619
$ws.Range("A6").Value = "2023-01-01 00:00:00_diff"
$ws.Range("B6").Value = -4.14234893673768
$ws.Range("C6").Value = -0.4417274786350349
$ws.Range("D6").Value = -0.6412421776795769
$ws.Range("E6").Value = 1.512882179556864
$ws.Range("F6").Value = -2.921454045212219
$ws.Range("G6").Value = 1.799337457552009
$ws.Range("H6").Value = 0.1376777883546749
$ws.Range("A7").Value = "2023-04-01 00:00:00_diff"
$ws.Range("B7").Value = -0.6288648973096044
$ws.Range("C7").Value = -0.8283795963541464
$ws.Range("D7").Value = 1.325744760882294
$ws.Range("E7").Value = -3.108591463886789
$ws.Range("F7").Value = 1.61220003887744
$ws.Range("G7").Value = -0.04945963031989459
$ws.Range("A8").Value = "2023-07-01 00:00:00_diff"
$ws.Range("B8").Value = -0.9420242164422319
$ws.Range("C8").Value = 1.212100140794209
$ws.Range("D8").Value = -3.222236083974874
$ws.Range("E8").Value = 1.498555418789354
$ws.Range("F8").Value = -0.1631042504079802
$ws.Range("G8").Value = -1.553285820388052
$ws.Range("H8").Value = -1.037158441845645
$ws.Range("I8").Value = 1.563884698290838
$ws.Range("A9").Value = "2023-10-01 00:00:00_diff"
$ws.Range("B9").Value = 0.9834580318344777
$ws.Range("C9").Value = -3.450878192934606
$ws.Range("D9").Value = 1.269913309829623
$ws.Range("E9").Value = -0.3917463593677113
$ws.Range("F9").Value = -1.781927929347783
$ws.Range("G9").Value = -1.265800550805376
$ws.Range("H9").Value = 1.335242589331107
$ws.Range("A10").Value = "2024-01-01 00:00:00_diff"
$ws.Range("B10").Value = -3.353876282308941
$ws.Range("C10").Value = 1.366915220455287
$ws.Range("D10").Value = -0.294744448742047
$ws.Range("E10").Value = -1.684926018722119
$ws.Range("F10").Value = -1.168798640179712
$ws.Range("G10").Value = 1.432244499956771
$ws.Range("A11").Value = "2024-04-01 00:00:00_diff"
$ws.Range("B11").Value = 1.609820403724652
$ws.Range("C11").Value = -0.05183926547268242
$ws.Range("D11").Value = -1.442020835452754
$ws.Range("E11").Value = -0.9258934569103474
$ws.Range("F11").Value = 1.675149683226136
$ws.Range("A12").Value = "2024-07-01 00:00:00_diff"
$ws.Range("B12").Value = -0.7294052785381211
$ws.Range("C12").Value = -2.119586848518193
$ws.Range("D12").Value = -1.603459469975786
$ws.Range("E12").Value = 0.997583670160697
$ws.Range("A13").Value = "2024-10-01 00:00:00_diff"
$ws.Range("B13").Value = -1.846431246597561
$ws.Range("C13").Value = -1.330303868055154
$ws.Range("D13").Value = 1.270739272081329
$ws.Range("A14").Value = "2025-01-01 00:00:00_diff"
$ws.Range("B14").Value = -1.121732059169287
$ws.Range("C14").Value = 1.479311080967196
$ws.Range("A15").Value = "2025-04-01 00:00:00_diff"
$ws.Range("B15").Value = 1.884604928957667
$ws.Range("A16").Value = "2025-07-01 00:00:00_diff"
